$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column R (modelName) - filled top to bottom first
$ws.Range("R1").Style = "Normal"
$ws.Range("R1").Value = "modelName"
$ws.Range("R2").Style = "Normal"
$ws.Range("R2").Value = "pseudobatracotoxin_molecule.xyz"
$ws.Range("R3").Style = "Normal"
$ws.Range("R3").Value = "pseudobatracotoxin_molecule.xyz"
$ws.Range("R4").Style = "Normal"
$ws.Range("R4").Value = "MRT_VK_mol.xyz"
$ws.Range("R5").Style = "Normal"
$ws.Range("R5").Value = "pseudobatracotoxin_molecule.xyz"
$ws.Range("R6").Style = "Normal"
$ws.Range("R6").Value = "pseudobatracotoxin_molecule.xyz"
$ws.Range("R7").Style = "Normal"
$ws.Range("R7").Value = "MRT_VK_mol.xyz"

# Header row for reference quaternion columns N:Q
$ws.Range("N1").Style = "Normal"
$ws.Range("N1").Value = "ref_i"
$ws.Range("O1").Style = "Normal"
$ws.Range("O1").Value = "ref_j"
$ws.Range("P1").Style = "Normal"
$ws.Range("P1").Value = "ref_k"
$ws.Range("Q1").Style = "Normal"
$ws.Range("Q1").Value = "ref_r"

# Reference quaternion data rows
$ws.Range("N2").Style = "Normal"
$ws.Range("N2").Value = -0.490257597533579
$ws.Range("O2").Style = "Normal"
$ws.Range("O2").Value = -0.28253363768972101
$ws.Range("P2").Style = "Normal"
$ws.Range("P2").Value = -0.45539483504230199
$ws.Range("Q2").Style = "Normal"
$ws.Range("Q2").Value = 0.68734109134490795

$ws.Range("N3").Style = "Normal"
$ws.Range("N3").Value = -0.490257597533579
$ws.Range("O3").Style = "Normal"
$ws.Range("O3").Value = -0.28253363768972101
$ws.Range("P3").Style = "Normal"
$ws.Range("P3").Value = -0.45539483504230199
$ws.Range("Q3").Style = "Normal"
$ws.Range("Q3").Value = 0.68734109134490795

$ws.Range("N4").Style = "Normal"
$ws.Range("N4").Value = 0.29922234471855402
$ws.Range("O4").Style = "Normal"
$ws.Range("O4").Value = 0.223654136694591
$ws.Range("P4").Style = "Normal"
$ws.Range("P4").Value = -0.10297736856098801
$ws.Range("Q4").Style = "Normal"
$ws.Range("Q4").Value = 0.92186792824394204

$ws.Range("N5").Style = "Normal"
$ws.Range("N5").Value = -0.490257597533579
$ws.Range("O5").Style = "Normal"
$ws.Range("O5").Value = -0.28253363768972101
$ws.Range("P5").Style = "Normal"
$ws.Range("P5").Value = -0.45539483504230199
$ws.Range("Q5").Style = "Normal"
$ws.Range("Q5").Value = 0.68734109134490795

$ws.Range("N6").Style = "Normal"
$ws.Range("N6").Value = -0.490257597533579
$ws.Range("O6").Style = "Normal"
$ws.Range("O6").Value = -0.28253363768972101
$ws.Range("P6").Style = "Normal"
$ws.Range("P6").Value = -0.45539483504230199
$ws.Range("Q6").Style = "Normal"
$ws.Range("Q6").Value = 0.68734109134490795

$ws.Range("N7").Style = "Normal"
$ws.Range("N7").Value = 0.29922234471855402
$ws.Range("O7").Style = "Normal"
$ws.Range("O7").Value = 0.223654136694591
$ws.Range("P7").Style = "Normal"
$ws.Range("P7").Value = -0.10297736856098801
$ws.Range("Q7").Style = "Normal"
$ws.Range("Q7").Value = 0.92186792824394204

# Update active cell selection as reflected in the saved view state
$ws.Range("A3").Select()
